# "finished week 1 analysis"
# Update exercise names on the Maintenance sheet to reflect week 1 changes,
# then move the active selection to D4 (matches author's final cursor spot).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value  = "Barbell Bench Press"      # was: Flat Bench Press
$ws.Range("A4").Value  = "Machine Fly"              # was: Chest Fly
$ws.Range("A5").Value  = "Incline Press"            # was: Incline Machine Press
$ws.Range("A7").Value  = "Cable Lateral Raise"      # was: Lateral Raise
$ws.Range("A8").Value  = "Cable Curl"               # was: Bicep Curl

$ws.Range("D15").Value = "Flat Press"               # was: Flat Machine Press

$ws.Range("A17").Value = "Face Pull"                # was: Rear Delt
$ws.Range("D17").Value = "Dumbbell Lateral Raise"   # was: Lateral Raise
$ws.Range("A18").Value = "Dumbbell Lateral Raise"   # was: Lateral Raise

$ws.Range("A23").Value = "Barbell Overhead Press"   # was: Shoulder Press
$ws.Range("A24").Value = "Pull Ups"                 # was: Pull ups
$ws.Range("A25").Value = "Cable Row"                # was: Seated Row
$ws.Range("A26").Value = "Dumbbell Curl"            # was: Bicep curl

$ws.Range("D4").Select()
